# "kebutuhan_kelas" used to keep the eligible-room list split into two
# separate strings (one for ordinary rooms, one for the CM/G5 lab rooms).
# This update merges them into a single combined list and uses it for
# every data row of column N ("rooms"), replacing whichever of the two
# partial lists (or stray value) a row used to have.
$wb = $excel.ActiveWorkbook
$wsDosen = $wb.Worksheets.Item("dosen")
$wsKelas = $wb.Worksheets.Item("kebutuhan_kelas")

$mergedRooms = "B2-R1, B3-R1, B3-R2, B3R3, CM-101, CM-102, CM-103, CM-201, CM-202, CM-203, CM-204, CM-205, CM-208, G2-R2, G2-R3, G2-R4, G2-R5, G2-R6, G2-R7, G3-R1, G3-R2, G3-R3, G3-R4, G4-R1, G4-R2, G4-R3, G4-R4, CM-206, CM-207, CM-LabVirtual, CM-Lab3, G5-Lab1, G5-Lab2"

# All data rows (2 through 374) of column N get the merged room list.
$wsKelas.Range("N2:N374").Value = $mergedRooms

# Row 209's N cell previously carried a stray "horizontal left" alignment
# that none of the surrounding cells have; restore the default/Normal
# style so it matches the rest of the column.
$wsKelas.Range("N209").Style = "Normal"

# "kebutuhan_kelas" becomes the active / selected sheet, with N2:N374
# selected, instead of "dosen" being active with a selection near B50.
$wsDosen.Range("F78").Select() | Out-Null
$wsKelas.Activate() | Out-Null
$wsKelas.Range("N2:N374").Select() | Out-Null
